$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 0.059691
$ws.Range("M2").Value = 2.027115333333333
$ws.Range("N2").Value = 6.081346
$ws.Range("O2").Value = 0.006596284565418616
$ws.Range("P2").Value = 0.006596284565418615
$ws.Range("Q2").Value = 0.04033351378733333
$ws.Range("R2").Value = 0.363001624086
$ws.Range("S2").Value = 0.006596284565418616
$ws.Range("T2").Value = 0.006596284565418615

# Row 3
$ws.Range("H3").Value = 0.059691
$ws.Range("O3").Value = 0.8344762556643375
$ws.Range("P3").Value = 0.8344762556643374
$ws.Range("Q3").Value = 5.102472343217333
$ws.Range("S3").Value = 0.8344762556643375
$ws.Range("T3").Value = 0.8344762556643374

# Row 4
$ws.Range("H4").Value = 0.059691
$ws.Range("M4").Value = 48.84026566666667
$ws.Range("N4").Value = 146.520797
$ws.Range("O4").Value = 0.158927459770244
$ws.Range("P4").Value = 0.158927459770244
$ws.Range("Q4").Value = 0.9717747659696669
$ws.Range("R4").Value = 8.745972893727002
$ws.Range("S4").Value = 0.158927459770244
$ws.Range("T4").Value = 0.158927459770244
